$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 6.5
$ws.Range("AP3").Value = 1.78
$ws.Range("AQ3").Value = 2.1
